$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new data columns (AC, AD) were appended: a header row value (new
# shared strings) plus the two probability values below each header.
$ws.Range("AC1").Value = "wnb-调节6Hz_20161230_113123_ASIC_EEG"
$ws.Range("AD1").Value = "wnb-调节6Hz_20170110_113300_ASIC_EEG"

$ws.Range("AC2").Value = 0.93890675241157551
$ws.Range("AD2").Value = 0.98381877022653719

$ws.Range("AC3").Value = 0.94409937888198758
$ws.Range("AD3").Value = 0.91808873720136519

# The sheet selection grows to cover the newly added columns.
$ws.Range("A1:AD3").Select() | Out-Null

# Best-effort: the saved workbook window size also changed in the source
# file (a side effect of the author's on-screen Excel window being resized
# before saving). Attempt to mirror it through the exposed window/app
# properties, in case the host honours it.
$excel.ActiveWindow.Width = 19170
$excel.ActiveWindow.Height = 11655
